$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2988413333333333
$ws.Range("H2").Value = 0.896524
$ws.Range("I2").Value = 0.3632971504731247
$ws.Range("J2").Value = 0.3632971504731246
$ws.Range("M2").Value = 0.003058333333333333
$ws.Range("N2").Value = 0.009175
$ws.Range("O2").Value = 0.0001379486413073712
$ws.Range("P2").Value = 0.0001379486413073712
$ws.Range("Q2").Value = 0.0009139564111111112
$ws.Range("R2").Value = 0.0082256077
$ws.Range("S2").Value = 0.00005011634829860715
$ws.Range("T2").Value = 0.00005011634829860714

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2988413333333333
$ws.Range("H3").Value = 0.896524
$ws.Range("I3").Value = 0.3632971504731247
$ws.Range("J3").Value = 0.3632971504731246
$ws.Range("O3").Value = 0.9939610820947024
$ws.Range("P3").Value = 0.9939610820947024
$ws.Range("Q3").Value = 6.585328385737778
$ws.Range("R3").Value = 59.26795547164
$ws.Range("S3").Value = 0.3611032288061889
$ws.Range("T3").Value = 0.3611032288061888

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2988413333333333
$ws.Range("H4").Value = 0.896524
$ws.Range("I4").Value = 0.3632971504731247
$ws.Range("J4").Value = 0.3632971504731246
$ws.Range("O4").Value = 0.005900969263990248
$ws.Range("P4").Value = 0.005900969263990248
$ws.Range("Q4").Value = 0.03909591743333334
$ws.Range("R4").Value = 0.3518632569
$ws.Range("S4").Value = 0.002143805318637149
$ws.Range("T4").Value = 0.002143805318637148

# Row 5
$ws.Range("I5").Value = 0.4682720202225272
$ws.Range("J5").Value = 0.4682720202225272
$ws.Range("M5").Value = 0.003058333333333333
$ws.Range("N5").Value = 0.009175
$ws.Range("O5").Value = 0.0001379486413073712
$ws.Range("P5").Value = 0.0001379486413073712
$ws.Range("Q5").Value = 0.001178044513888889
$ws.Range("R5").Value = 0.010602400625
$ws.Range("S5").Value = 0.00006459748895195549
$ws.Range("T5").Value = 0.00006459748895195549

# Row 6
$ws.Range("I6").Value = 0.4682720202225272
$ws.Range("J6").Value = 0.4682720202225272
$ws.Range("O6").Value = 0.9939610820947024
$ws.Range("P6").Value = 0.9939610820947024
$ws.Range("S6").Value = 0.4654441639350555
$ws.Range("T6").Value = 0.4654441639350555

# Row 7
$ws.Range("I7").Value = 0.4682720202225272
$ws.Range("J7").Value = 0.4682720202225272
$ws.Range("O7").Value = 0.005900969263990248
$ws.Range("P7").Value = 0.005900969263990248
$ws.Range("S7").Value = 0.002763258798519753
$ws.Range("T7").Value = 0.002763258798519753

# Row 8
$ws.Range("I8").Value = 0.1684308293043481
$ws.Range("J8").Value = 0.1684308293043481
$ws.Range("M8").Value = 0.003058333333333333
$ws.Range("N8").Value = 0.009175
$ws.Range("O8").Value = 0.0001379486413073712
$ws.Range("P8").Value = 0.0001379486413073712
$ws.Range("Q8").Value = 0.0004237259666666667
$ws.Range("R8").Value = 0.0038135337
$ws.Range("S8").Value = 0.00002323480405680859
$ws.Range("T8").Value = 0.00002323480405680859

# Row 9
$ws.Range("I9").Value = 0.1684308293043481
$ws.Range("J9").Value = 0.1684308293043481
$ws.Range("O9").Value = 0.9939610820947024
$ws.Range("P9").Value = 0.9939610820947024
$ws.Range("Q9").Value = 3.053071899426667
$ws.Range("S9").Value = 0.167413689353458
$ws.Range("T9").Value = 0.167413689353458

# Row 10
$ws.Range("I10").Value = 0.1684308293043481
$ws.Range("J10").Value = 0.1684308293043481
$ws.Range("O10").Value = 0.005900969263990248
$ws.Range("P10").Value = 0.005900969263990248
$ws.Range("S10").Value = 0.0009939051468333463
$ws.Range("T10").Value = 0.0009939051468333463
